$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the date/time value in row 89 (column A) - truncate to just the date portion
$ws.Range("A89").Value = 45464.2916666667

# Append a new row 90 with the latest data point
$ws.Range("A90").Value = 45467.6494791667
$ws.Range("B90").Value = 14100
$ws.Range("C90").Value = 6.5
$ws.Range("D90").Value = 6.26000022888184
$ws.Range("E90").Value = 6.07999992370605
$ws.Range("F90").Value = 6.28000020980835
$ws.Range("G90").Value = "'6.28000020980835"
$ws.Range("G90").Style = "Normal"
$ws.Range("H90").Value = "PAL.MI"

# Copy the date-format style from A89 into the new A90 cell
$ws.Range("A89").Copy()
$ws.Range("A90").PasteSpecial(-4122)
$ws.Range("A90").Value = 45467.6494791667
